$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update cell B18 from "%" to ">"
$ws.Range("B18").Value = ">"
